# Update the Taxable Account dividend amount for January 2017 on the
# "Yearly" sheet. All dependent totals (Grand Total for that row, the
# yearly sum row, and the "All Time" sheet which references Yearly!L3)
# recalculate automatically from this single source value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Yearly")
$ws.Range("L3").Value = 57.62

$excel.CalculateFullRebuild()
